# Update the "想去人数" (want-to-go count) figures in column F for the
# exhibition data, which is duplicated across the "展览" and "全部类型"
# worksheets.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    5  = 247
    7  = 112
    13 = 80
    14 = 337
    15 = 39
    16 = 464
    17 = 385
    19 = 59
    20 = 30
    22 = 919
    23 = 2707
    26 = 520
    27 = 962
    30 = 257
    33 = 592
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
